# Update the multiplication-problem table: each populated cell's
# "NNN×N=" expression is replaced with a new one, per the commit's
# regenerated numbers. Rows 1, 5, 10, 15, 20 (1-based) hold the five
# problems each; the other rows are blank spacer rows.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "959×8="; New = "848×6=" },
    @{ Row = 1;  Col = 2; Old = "234×5="; New = "157×3=" },
    @{ Row = 1;  Col = 3; Old = "301×2="; New = "155×3=" },
    @{ Row = 1;  Col = 4; Old = "809×8="; New = "204×6=" },
    @{ Row = 1;  Col = 5; Old = "795×6="; New = "611×4=" },

    @{ Row = 5;  Col = 1; Old = "234×5="; New = "682×3=" },
    @{ Row = 5;  Col = 2; Old = "580×4="; New = "483×3=" },
    @{ Row = 5;  Col = 3; Old = "280×2="; New = "478×4=" },
    @{ Row = 5;  Col = 4; Old = "260×9="; New = "234×9=" },
    @{ Row = 5;  Col = 5; Old = "931×7="; New = "746×3=" },

    @{ Row = 10; Col = 1; Old = "944×4="; New = "443×3=" },
    @{ Row = 10; Col = 2; Old = "569×4="; New = "106×2=" },
    @{ Row = 10; Col = 3; Old = "684×8="; New = "548×6=" },
    @{ Row = 10; Col = 4; Old = "761×2="; New = "301×6=" },
    @{ Row = 10; Col = 5; Old = "406×5="; New = "424×2=" },

    @{ Row = 15; Col = 1; Old = "727×4="; New = "445×8=" },
    @{ Row = 15; Col = 2; Old = "186×6="; New = "370×2=" },
    @{ Row = 15; Col = 3; Old = "764×2="; New = "189×5=" },
    @{ Row = 15; Col = 4; Old = "687×9="; New = "906×7=" },
    @{ Row = 15; Col = 5; Old = "548×2="; New = "762×8=" },

    @{ Row = 20; Col = 1; Old = "870×4="; New = "387×7=" },
    @{ Row = 20; Col = 2; Old = "108×7="; New = "843×3=" },
    @{ Row = 20; Col = 3; Old = "473×3="; New = "795×4=" },
    @{ Row = 20; Col = 4; Old = "877×8="; New = "158×7=" },
    @{ Row = 20; Col = 5; Old = "513×7="; New = "665×3=" }
)

foreach ($chg in $changes) {
    $cell = $t.Cell($chg.Row, $chg.Col)
    $r = $cell.Range
    $r.Text = $chg.New
}
